$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (row 7) with the same shape/style as the existing rows.
$row = 7

$ws.Cells.Item($row, 1).Value = 42612.895370370374
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item($row, 2).Value = -8
$ws.Cells.Item($row, 3).Value = 51
$ws.Cells.Item($row, 4).Value = 47
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 100
$ws.Cells.Item($row, 7).Value = 19474
$ws.Cells.Item($row, 8).Value = 10212
$ws.Cells.Item($row, 9).Value = 544
$ws.Cells.Item($row, 10).Value = 105
$ws.Cells.Item($row, 11).Value = 97
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 4
$ws.Cells.Item($row, 14).Value = "Named"
